$d = $word.ActiveDocument

# Locate the paragraph that ends with "...not just lucky run."
# (3rd paragraph in the document) and insert a brand new bullet
# paragraph right after it, inheriting the same list/paragraph
# formatting (NoSpacing style, bullet numId=1, bold+underline rPr).
$sourcePara = $d.Paragraphs(3)
$insertionPoint = $d.Range($sourcePara.Range.End - 1, $sourcePara.Range.End - 1)
$insertionPoint.InsertParagraphAfter()

# The newly created paragraph is now paragraph 4; fill in its text.
$newPara = $d.Paragraphs(4)
$newParaText = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$newParaText.InsertAfter("Prepared presentation for invited lecture on 19")

$newPara2 = $d.Paragraphs(4)
$newParaEnd = $d.Range($newPara2.Range.End - 1, $newPara2.Range.End - 1)
$newParaEnd.InsertAfter("/03.")
